# Update the "Last Updated" date shown on every slide from 20-Apr-20 to
# 21-Apr-20, without disturbing any other runs/formatting in the textbox
# (the shape also contains the project name/manager text and a hyperlink
# run, so we can't just overwrite TextFrame.TextRange.Text wholesale).

$p = $ppt.ActivePresentation

$oldDate = "20-Apr-20"
$newDate = "21-Apr-20"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $textFrame = $shape.TextFrame
        if (-not $textFrame.HasText) {
            continue
        }

        $textRange = $textFrame.TextRange
        $fullText = $textRange.Text

        if ([string]::IsNullOrEmpty($fullText)) {
            continue
        }

        $idx = $fullText.IndexOf($oldDate)
        while ($idx -ge 0) {
            $part = $textRange.Characters($idx + 1, $oldDate.Length)
            $part.Text = $newDate

            $fullText = $textRange.Text
            if ([string]::IsNullOrEmpty($fullText)) {
                $idx = -1
            } else {
                $idx = $fullText.IndexOf($oldDate)
            }
        }
    }
}
